# Allocation rule summary tables update: add "Within 5 miles" and
# "Within 10 miles of HFC production facility" columns (F and G) to both
# the "Means" and "Standard Deviations" worksheets, and refresh the
# existing values that changed with the new radius data included.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Means"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Means")

# Header row
$ws1.Cells.Item(1, 6).Value = "Within 5 miles of HFC production facility"
$ws1.Cells.Item(1, 7).Value = "Within 10 miles of HFC production facility"

# % White
$ws1.Cells.Item(2, 6).Value = 90
$ws1.Cells.Item(2, 7).Value = 81

# % Black or African American
$ws1.Cells.Item(3, 6).Value = 4.3
$ws1.Cells.Item(3, 7).Value = 13

# % Other
$ws1.Cells.Item(4, 6).Value = 5.2
$ws1.Cells.Item(4, 7).Value = 5.8

# % Hispanic
$ws1.Cells.Item(5, 6).Value = 3.4
$ws1.Cells.Item(5, 7).Value = 3.3

# Median Income [1,000 2019$]
$ws1.Cells.Item(6, 6).Value = 71
$ws1.Cells.Item(6, 7).Value = 66

# % Below Poverty Line
$ws1.Cells.Item(7, 6).Value = 5.5
$ws1.Cells.Item(7, 7).Value = 7.5

# % Below Half the Poverty Line
$ws1.Cells.Item(8, 6).Value = 5.5
$ws1.Cells.Item(8, 7).Value = 7.6

# Total Cancer Risk (per million) - existing cols B & C also updated
$ws1.Cells.Item(9, 2).Value = 29
$ws1.Cells.Item(9, 3).Value = 29
$ws1.Cells.Item(9, 6).Value = 30
$ws1.Cells.Item(9, 7).Value = 30

# Total Respiratory (hazard quotient) - existing cols B-E also updated
$ws1.Cells.Item(10, 2).Value = 0.37
$ws1.Cells.Item(10, 3).Value = 0.36
$ws1.Cells.Item(10, 4).Value = 0.38
$ws1.Cells.Item(10, 5).Value = 0.38
$ws1.Cells.Item(10, 6).Value = 0.39
$ws1.Cells.Item(10, 7).Value = 0.38

# ---------------------------------------------------------------------
# Sheet "Standard Deviations"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# Header row
$ws2.Cells.Item(1, 6).Value = "Within 5 mile of HFC production facility SD"
$ws2.Cells.Item(1, 7).Value = "Within 10 mile of HFC production facility SD"

# % White
$ws2.Cells.Item(2, 6).Value = 12
$ws2.Cells.Item(2, 7).Value = 28

# % Black or African American
$ws2.Cells.Item(3, 6).Value = 7.4
$ws2.Cells.Item(3, 7).Value = 27

# % Other
$ws2.Cells.Item(4, 6).Value = 8.3
$ws2.Cells.Item(4, 7).Value = 6.9

# % Hispanic
$ws2.Cells.Item(5, 6).Value = 6.7
$ws2.Cells.Item(5, 7).Value = 6.4

# Median Income [1,000 2019$]
$ws2.Cells.Item(6, 6).Value = 26
$ws2.Cells.Item(6, 7).Value = 31

# % Below Poverty Line
$ws2.Cells.Item(7, 6).Value = 8
$ws2.Cells.Item(7, 7).Value = 11

# % Below Half the Poverty Line
$ws2.Cells.Item(8, 6).Value = 6.9
$ws2.Cells.Item(8, 7).Value = 11

# Total Cancer Risk (per million) - existing cols B & D & E also updated
$ws2.Cells.Item(9, 2).Value = 10
$ws2.Cells.Item(9, 4).Value = 0
$ws2.Cells.Item(9, 5).Value = 0
$ws2.Cells.Item(9, 6).Value = 0
$ws2.Cells.Item(9, 7).Value = 1.9

# Total Respiratory (hazard quotient) - existing cols B-E also updated
$ws2.Cells.Item(10, 2).Value = 0.14
$ws2.Cells.Item(10, 3).Value = 0.086
$ws2.Cells.Item(10, 4).Value = 0.041
$ws2.Cells.Item(10, 5).Value = 0.037
$ws2.Cells.Item(10, 6).Value = 0.039
$ws2.Cells.Item(10, 7).Value = 0.039
